$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-04-08 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-04-09 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("65×45=2925", $true, $false, $false, $false, $false, $true, 1, $false, "23×32=736", 2) | Out-Null
$d.Content.Find.Execute("75×65=4875", $true, $false, $false, $false, $false, $true, 1, $false, "19×67=1273", 2) | Out-Null
$d.Content.Find.Execute("31×86=2666", $true, $false, $false, $false, $false, $true, 1, $false, "99×35=3465", 2) | Out-Null
$d.Content.Find.Execute("57×82=4674", $true, $false, $false, $false, $false, $true, 1, $false, "62×37=2294", 2) | Out-Null
$d.Content.Find.Execute("98×45=4410", $true, $false, $false, $false, $false, $true, 1, $false, "25×32=800", 2) | Out-Null
$d.Content.Find.Execute("95×73=6935", $true, $false, $false, $false, $false, $true, 1, $false, "91×12=1092", 2) | Out-Null
$d.Content.Find.Execute("64×62=3968", $true, $false, $false, $false, $false, $true, 1, $false, "57×88=5016", 2) | Out-Null
$d.Content.Find.Execute("60×98=5880", $true, $false, $false, $false, $false, $true, 1, $false, "86×72=6192", 2) | Out-Null
$d.Content.Find.Execute("12×91=1092", $true, $false, $false, $false, $false, $true, 1, $false, "95×87=8265", 2) | Out-Null
$d.Content.Find.Execute("91×81=7371", $true, $false, $false, $false, $false, $true, 1, $false, "53×20=1060", 2) | Out-Null
$d.Content.Find.Execute("42×37=1554", $true, $false, $false, $false, $false, $true, 1, $false, "93×82=7626", 2) | Out-Null
$d.Content.Find.Execute("50×89=4450", $true, $false, $false, $false, $false, $true, 1, $false, "68×26=1768", 2) | Out-Null
$d.Content.Find.Execute("65×71=4615", $true, $false, $false, $false, $false, $true, 1, $false, "48×65=3120", 2) | Out-Null
$d.Content.Find.Execute("75×35=2625", $true, $false, $false, $false, $false, $true, 1, $false, "98×20=1960", 2) | Out-Null
$d.Content.Find.Execute("21×48=1008", $true, $false, $false, $false, $false, $true, 1, $false, "35×47=1645", 2) | Out-Null
$d.Content.Find.Execute("50×80=4000", $true, $false, $false, $false, $false, $true, 1, $false, "55×16=880", 2) | Out-Null
$d.Content.Find.Execute("63×18=1134", $true, $false, $false, $false, $false, $true, 1, $false, "66×61=4026", 2) | Out-Null
$d.Content.Find.Execute("58×16=928", $true, $false, $false, $false, $false, $true, 1, $false, "86×82=7052", 2) | Out-Null
$d.Content.Find.Execute("45×83=3735", $true, $false, $false, $false, $false, $true, 1, $false, "74×42=3108", 2) | Out-Null
$d.Content.Find.Execute("72×65=4680", $true, $false, $false, $false, $false, $true, 1, $false, "59×36=2124", 2) | Out-Null
$d.Content.Find.Execute("47×42=1974", $true, $false, $false, $false, $false, $true, 1, $false, "64×86=5504", 2) | Out-Null
$d.Content.Find.Execute("64×24=1536", $true, $false, $false, $false, $false, $true, 1, $false, "20×55=1100", 2) | Out-Null
$d.Content.Find.Execute("18×92=1656", $true, $false, $false, $false, $false, $true, 1, $false, "26×73=1898", 2) | Out-Null
$d.Content.Find.Execute("91×60=5460", $true, $false, $false, $false, $false, $true, 1, $false, "95×86=8170", 2) | Out-Null
$d.Content.Find.Execute("26×65=1690", $true, $false, $false, $false, $false, $true, 1, $false, "60×51=3060", 2) | Out-Null
